$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-02 Tuesday" "2024-01-03 Wednesday"

Replace-Text "99×38=3762" "54×34=1836"
Replace-Text "81×70=5670" "97×20=1940"
Replace-Text "38×37=1406" "51×48=2448"
Replace-Text "79×19=1501" "68×49=3332"
Replace-Text "14×75=1050" "57×75=4275"
Replace-Text "58×65=3770" "72×53=3816"
Replace-Text "35×33=1155" "69×76=5244"
Replace-Text "77×73=5621" "32×66=2112"
Replace-Text "69×83=5727" "67×31=2077"
Replace-Text "34×94=3196" "48×42=2016"
Replace-Text "71×61=4331" "16×49=784"
Replace-Text "33×50=1650" "21×85=1785"
Replace-Text "78×88=6864" "84×44=3696"
Replace-Text "64×22=1408" "92×54=4968"
Replace-Text "43×52=2236" "56×18=1008"
Replace-Text "72×44=3168" "60×55=3300"
Replace-Text "67×56=3752" "16×54=864"
Replace-Text "97×71=6887" "31×84=2604"
Replace-Text "24×34=816" "79×95=7505"
Replace-Text "17×42=714" "14×51=714"
Replace-Text "79×52=4108" "19×90=1710"
Replace-Text "85×65=5525" "31×71=2201"
Replace-Text "22×61=1342" "53×35=1855"
Replace-Text "95×88=8360" "39×62=2418"
Replace-Text "82×60=4920" "35×13=455"
